$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.091562213532490944,
    -0.0059999999657414094,
    -0.0039999999732387437,
    -0.0079999999498117091,
    -0.0029999999769492192,
    -0.0019999999811552982,
    -0.009999999937133186,
    -0.0099999999389437377,
    -0.001999999987794876,
    0.059518927100990382,
    -0.0029999999870122807,
    -0.0034999999836147211,
    -0.0034999999815470417,
    -0.0079999999559614565,
    0.003216645502210902,
    -0.0019999999889059872,
    -0.0019999999883237862,
    -0.0039999999769202432,
    -0.0039999999769158023,
    -0.0039999999773048245,
    0.012707284690893772,
    -0.0039999999773687733,
    -0.00085977482048882337,
    -0.019999999878793417,
    -0.019999999877073904,
    -0.0024999999777737258,
    -0.0024999999769907966,
    -0.0019999999767632559,
    -0.0069999999466370255,
    -0.059999999649423597,
    -0.0069999999496239695,
    -0.0099999999342639256,
    0.055639658857961294
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
